# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact string representation (avoid Excel
# auto-converting numeric-looking text into real numbers).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") "27.938.73"
Set-TextValue $ws.Range("E2") "  -1.88%  "
Set-TextValue $ws.Range("D3") "1.766.22"
Set-TextValue $ws.Range("E3") "  -3.44%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "321.41"
Set-TextValue $ws.Range("E5") "  -2.70%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.19%  "
Set-TextValue $ws.Range("D7") "0.4280"
Set-TextValue $ws.Range("E7") "  -4.68%  "
Set-TextValue $ws.Range("D8") "0.3613"
Set-TextValue $ws.Range("E8") "  -4.85%  "
Set-TextValue $ws.Range("D9") "43.20"
Set-TextValue $ws.Range("E9") "  -3.90%  "
Set-TextValue $ws.Range("D10") "0.07450"
Set-TextValue $ws.Range("E10") "  -4.65%  "
Set-TextValue $ws.Range("D11") "1.093"
Set-TextValue $ws.Range("E11") "  -4.31%  "
Set-TextValue $ws.Range("D12") "1.003"
Set-TextValue $ws.Range("E12") "  +0.28%  "
Set-TextValue $ws.Range("D13") "21.06"
Set-TextValue $ws.Range("E13") "  -5.80%  "
Set-TextValue $ws.Range("D14") "6.087"
Set-TextValue $ws.Range("E14") "  -4.91%  "
Set-TextValue $ws.Range("D15") "7.313"
Set-TextValue $ws.Range("E15") "  -3.53%  "
Set-TextValue $ws.Range("D16") "1.796.11"
Set-TextValue $ws.Range("E16") "  -2.30%  "
Set-TextValue $ws.Range("D17") "91.84"
Set-TextValue $ws.Range("E17") "  -2.49%  "
Set-TextValue $ws.Range("D18") "0.00001059"
Set-TextValue $ws.Range("E18") "  -2.86%  "
Set-TextValue $ws.Range("E19") "  +0.39%  "
Set-TextValue $ws.Range("E20") "  -0.01%  "
Set-TextValue $ws.Range("D21") "17.14"
Set-TextValue $ws.Range("E21") "  -2.88%  "
Set-TextValue $ws.Range("D22") "5.957"
Set-TextValue $ws.Range("E22") "  -7.14%  "
Set-TextValue $ws.Range("D23") "27.959.78"
Set-TextValue $ws.Range("E23") "  -2.05%  "
Set-TextValue $ws.Range("D24") "11.30"
Set-TextValue $ws.Range("E24") "  -4.52%  "
Set-TextValue $ws.Range("D25") "2.101"
Set-TextValue $ws.Range("E25") "  -7.93%  "
Set-TextValue $ws.Range("D26") "157.87"
Set-TextValue $ws.Range("E26") "  +2.37%  "
Set-TextValue $ws.Range("D27") "20.21"
Set-TextValue $ws.Range("E27") "  -3.73%  "
Set-TextValue $ws.Range("D28") "1.990.09"
Set-TextValue $ws.Range("E28") "  -2.67%  "
Set-TextValue $ws.Range("E29") "  -9.52%  "
Set-TextValue $ws.Range("D30") "125.39"
Set-TextValue $ws.Range("E30") "  -3.40%  "
Set-TextValue $ws.Range("D31") "1.158"
Set-TextValue $ws.Range("E31") "  -4.07%  "
Set-TextValue $ws.Range("D32") "3.747"
Set-TextValue $ws.Range("E32") "  +2.06%  "
Set-TextValue $ws.Range("D33") "5.606"
Set-TextValue $ws.Range("E33") "  -4.96%  "
Set-TextValue $ws.Range("D34") "0.08878"
Set-TextValue $ws.Range("E34") "  -4.80%  "
Set-TextValue $ws.Range("D35") "12.47"
Set-TextValue $ws.Range("E35") "  -3.46%  "
Set-TextValue $ws.Range("D36") "0.02308"
Set-TextValue $ws.Range("E36") "  -2.56%  "
Set-TextValue $ws.Range("D37") "0.2105"
Set-TextValue $ws.Range("E37") "  -4.92%  "
Set-TextValue $ws.Range("D41") "1.183"
Set-TextValue $ws.Range("E41") "  -1.44%  "
Set-TextValue $ws.Range("D42") "1.001"
Set-TextValue $ws.Range("E42") "  +0.19%  "
Set-TextValue $ws.Range("E43") "  -0.54%  "
Set-TextValue $ws.Range("D44") "7.801"
Set-TextValue $ws.Range("E44") "  -4.71%  "
Set-TextValue $ws.Range("D45") "13.42"
Set-TextValue $ws.Range("E45") "  -4.32%  "
Set-TextValue $ws.Range("D46") "0.5920"
Set-TextValue $ws.Range("E46") "  -4.25%  "
Set-TextValue $ws.Range("D47") "3.696"
Set-TextValue $ws.Range("E47") "  -2.30%  "
Set-TextValue $ws.Range("D48") "2.002"
Set-TextValue $ws.Range("E48") "  -2.68%  "
Set-TextValue $ws.Range("D49") "122.42"
Set-TextValue $ws.Range("E49") "  -4.56%  "
Set-TextValue $ws.Range("E50") "  +2.40%  "
Set-TextValue $ws.Range("E51") "  -2.29%  "

# Rows 38-40: coin ordering changed (values rotated among the three rows)
Set-TextValue $ws.Range("B38") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "5.015"
Set-TextValue $ws.Range("E38") "  -4.41%  "
Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.06016"
Set-TextValue $ws.Range("E39") "  -4.81%  "
Set-TextValue $ws.Range("B40") "TheSandbox"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.6351"
Set-TextValue $ws.Range("E40") "  -5.31%  "
